$p = $ppt.ActivePresentation

# --- 1) Swap the active (slide-facing) theme's colour scheme from
#        "Red Violet"/Integral to the plain "Office" palette.
#        The font scheme / format scheme are already identical between
#        the two theme parts, so only the 12 scheme colours change.
$tcs = $p.Slides.Item(1).ThemeColorScheme
$officeColors = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)
for ($i = 1; $i -le 12; $i++) {
    $tcs.Colors($i).RGB = $officeColors[$i - 1]
}

# --- 2) Re-style the three tables (slides 14-16) from the default
#        table style to the new built-in style.
$newStyleId = "{C548CE67-6A3D-4804-BE37-375D0EAB3277}"
foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    foreach ($shape in $slide.Shapes) {
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}
